# Locate the shape whose text is exactly "RESTful App for Brand".
# It lives two levels deep (top-level group -> grouped rectangle), and the
# COM emulation here flattens nested sub-groups into a single GroupItems
# level, so a simple two-level iterative scan (no recursion, to keep this
# fast) is enough to find it regardless of shape ordering.
function Find-TextShape($topShapes, $searchText) {
    for ($i = 1; $i -le $topShapes.Count; $i++) {
        $shp = $topShapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and ($shp.TextFrame.TextRange.Text -eq $searchText)) {
            return $shp
        }
        if ($shp.Type -eq 6) {
            $gi = $shp.GroupItems
            for ($j = 1; $j -le $gi.Count; $j++) {
                $sub = $gi.Item($j)
                if ($sub.HasTextFrame -and $sub.TextFrame.HasText -and ($sub.TextFrame.TextRange.Text -eq $searchText)) {
                    return $sub
                }
            }
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$oldText = "RESTful App for Brand"
$keepPrefix = "RESTful App "
$newSuffix = "for Brand devices"

$shp = Find-TextShape $s.Shapes $oldText
if ($shp -eq $null) {
    # Fallback to the known location if the search above ever fails to match
    # (e.g. due to unexpected whitespace normalization).
    $shp = $s.Shapes.Item(7).GroupItems.Item("Rectangle 55")
}

$tr = $shp.TextFrame.TextRange
$fullLen = $tr.Text.Length
$prefixLen = $keepPrefix.Length

# Split the single run into two runs: keep "RESTful App " as-is (first run,
# unchanged formatting) and replace the remaining characters ("for Brand")
# with the new run text "for Brand devices" (second run), matching the
# authored edit which split one <a:r> into two.
$suffixRange = $tr.Characters($prefixLen + 1, $fullLen - $prefixLen)
$suffixRange.Text = $newSuffix
